$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.181.13'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").Value = '1.893.93'

$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.59'
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.681'
$ws.Range("E6").Value = '  +7.66%  '

$ws.Range("E7").Value = '  -0.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.71'
$ws.Range("E8").Value = '  -4.02%  '

$ws.Range("E9").Value = '  +2.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.86'
$ws.Range("E10").Value = '  +10.77%  '

$ws.Range("E11").Value = '  +1.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0981'
$ws.Range("E12").Value = '  -1.62%  '

$ws.Range("D13").Value = '2.170.07'
$ws.Range("E13").Value = '  -0.66%  '

$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.700'
$ws.Range("E15").Value = '  +1.26%  '

$ws.Range("D16").Value = '1.895.03'
$ws.Range("E16").Value = '  -1.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.77'
$ws.Range("E17").Value = '  -1.75%  '

$ws.Range("D18").Value = '35.185.47'
$ws.Range("E18").Value = '  -1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.80'
$ws.Range("E19").Value = '  -0.26%  '

$ws.Range("E20").Value = '  +0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '239.79'
$ws.Range("E21").Value = '  -1.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.54'
$ws.Range("E22").Value = '  +0.75%  '

$ws.Range("E23").Value = '  -3.50%  '

$ws.Range("E24").Value = '  -0.28%  '

$ws.Range("E25").Value = '  +0.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.31'
$ws.Range("E26").Value = '  +4.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.56'
$ws.Range("E27").Value = '  -2.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.52'
$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("E29").Value = '  +1.30%  '

$ws.Range("E30").Value = '  +3.39%  '

$ws.Range("E31").Value = '  +20.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.13'
$ws.Range("E32").Value = '  +0.68%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("E34").Value = '  +8.05%  '

$ws.Range("E35").Value = '  -0.30%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.916'
$ws.Range("E36").Value = '  -6.20%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.07'
$ws.Range("E37").Value = '  -1.99%  '

$ws.Range("E38").Value = '  +12.77%  '

$ws.Range("E39").Value = '  -0.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.42'
$ws.Range("E40").Value = '  +6.54%  '

$ws.Range("E41").Value = '  -2.06%  '

$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0630'
$ws.Range("E43").Value = '  +6.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '89.44'
$ws.Range("E44").Value = '  -2.24%  '

$ws.Range("D45").Value = '1.348.78'
$ws.Range("E45").Value = '  -0.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.40'
$ws.Range("E46").Value = '  +2.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.78'
$ws.Range("E48").Value = '  +0.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.63'
$ws.Range("E49").Value = '  -15.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '12.08'
$ws.Range("E50").Value = '  -4.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.43'
$ws.Range("E51").Value = '  -3.29%  '
